$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21: update dct:modified timestamp ---
$ws.Range('B21').Value = '2023-09-13T16:09:22+00:00'

# --- Row 24: subject -> analytical method ---
$ws.Range('B24').Value = 'analytical method'

# --- Row 25: variable -> identification method; drop notation 'var'; add broader vocab:1000 ---
$ws.Range('B25').Value = 'identification method'
$ws.Range('D25').ClearContents()
$ws.Range('G25').Value = 'vocab:1000'

# --- Row 26: community maturity level -> spectrometry; add broader vocab:1001 ---
$ws.Range('B26').Value = 'spectrometry'
$ws.Range('G26').Value = 'vocab:1001'

# --- Row 27: emerging -> mass spectrometry; drop altLabel 'new' (G27 broader vocab:1002 unchanged) ---
$ws.Range('B27').Value = 'mass spectrometry'
$ws.Range('C27').ClearContents()

# --- Row 28: developing -> inductively coupled plasma mass spectrometry; drop altLabel 'intermediate'; broader vocab:1002 -> vocab:1003 ---
$ws.Range('B28').Value = 'inductively coupled plasma mass spectrometry'
$ws.Range('C28').ClearContents()
$ws.Range('G28').Value = 'vocab:1003'

# --- Row 29: mature -> atmospheric presure chemical ionization high resolution mass spectrometry; broader vocab:1002 -> vocab:1003 ---
$ws.Range('B29').Value = 'atmospheric presure chemical ionization high resolution mass spectrometry'
$ws.Range('G29').Value = 'vocab:1003'

# --- Row 30: hasMaturityLevel -> high resolution mass spectrometry; add broader vocab:1003 ---
$ws.Range('B30').Value = 'high resolution mass spectrometry'
$ws.Range('G30').Value = 'vocab:1003'

# --- New rows 31-59: additional analytical-method vocabulary terms ---
$ws.Range('A31').Value = 'vocab:1007'
$ws.Range('B31').Value = 'electrothermal capture negative ion mass spectrometry'
$ws.Range('G31').Value = 'vocab:1003'
$ws.Range('A32').Value = 'vocab:1008'
$ws.Range('B32').Value = 'tandem mass spectrometry'
$ws.Range('G32').Value = 'vocab:1003'
$ws.Range('A33').Value = 'vocab:1009'
$ws.Range('B33').Value = 'mass spectrometry with Negative chemical ionization'
$ws.Range('G33').Value = 'vocab:1003'
$ws.Range('A34').Value = 'vocab:1010'
$ws.Range('B34').Value = 'other spectrometry '
$ws.Range('G34').Value = 'vocab:1001'
$ws.Range('A35').Value = 'vocab:1011'
$ws.Range('B35').Value = 'Flame Atomic Emission Spectrometry'
$ws.Range('G35').Value = 'vocab:1011'
$ws.Range('A36').Value = 'vocab:1012'
$ws.Range('B36').Value = 'flame Atomic absorption Spectrometry'
$ws.Range('G36').Value = 'vocab:1011'
$ws.Range('A37').Value = 'vocab:1013'
$ws.Range('B37').Value = 'electrotermal atomic absorption spectrometry'
$ws.Range('G37').Value = 'vocab:1011'
$ws.Range('A38').Value = 'vocab:1014'
$ws.Range('B38').Value = 'atomic absorption spectrometry'
$ws.Range('G38').Value = 'vocab:1011'
$ws.Range('A39').Value = 'vocab:1015'
$ws.Range('B39').Value = 'inductively coupled plasma with optical emission spectrometry'
$ws.Range('G39').Value = 'vocab:1011'
$ws.Range('A40').Value = 'vocab:1016'
$ws.Range('B40').Value = 'separation method'
$ws.Range('G40').Value = 'vocab:1000'
$ws.Range('A41').Value = 'vocab:1017'
$ws.Range('B41').Value = 'chromatography'
$ws.Range('G41').Value = 'vocab:1016'
$ws.Range('A42').Value = 'vocab:1018'
$ws.Range('B42').Value = 'gas chromatography'
$ws.Range('G42').Value = 'vocab:1017'
$ws.Range('A43').Value = 'vocab:1019'
$ws.Range('B43').Value = 'high performance liquid chromatography'
$ws.Range('G43').Value = 'vocab:1017'
$ws.Range('A44').Value = 'vocab:1020'
$ws.Range('B44').Value = 'combined method'
$ws.Range('G44').Value = 'vocab:1000'
$ws.Range('A45').Value = 'vocab:1021'
$ws.Range('B45').Value = 'HPLC combination'
$ws.Range('G45').Value = 'vocab:1020'
$ws.Range('A46').Value = 'vocab:1022'
$ws.Range('B46').Value = 'high performance liquid chromatography with mass spectrometry'
$ws.Range('G46').Value = 'vocab:1021'
$ws.Range('A47').Value = 'vocab:1023'
$ws.Range('B47').Value = 'high performance liquid chromatography with diode array detector'
$ws.Range('G47').Value = 'vocab:1021'
$ws.Range('A48').Value = 'vocab:1024'
$ws.Range('B48').Value = 'high performance liquid chromatography with high resolution mass spectrometry'
$ws.Range('G48').Value = 'vocab:1022'
$ws.Range('A49').Value = 'vocab:1025'
$ws.Range('B49').Value = 'high performance liquid chromatography with fluorescence detection'
$ws.Range('G49').Value = 'vocab:1021'
$ws.Range('A50').Value = 'vocab:1026'
$ws.Range('B50').Value = 'high performance liquid chromatography with tandem mass spectrometry'
$ws.Range('G50').Value = 'vocab:1022'
$ws.Range('A51').Value = 'vocab:1027'
$ws.Range('B51').Value = 'GC combination'
$ws.Range('G51').Value = 'vocab:1020'
$ws.Range('A52').Value = 'vocab:1028'
$ws.Range('B52').Value = 'gas chromatography with mass spectrometry'
$ws.Range('G52').Value = 'vocab:1027'
$ws.Range('A53').Value = 'vocab:1029'
$ws.Range('B53').Value = 'gas chromatography with electron capture detector'
$ws.Range('G53').Value = 'vocab:1028'
$ws.Range('A54').Value = 'vocab:1030'
$ws.Range('B54').Value = 'gas chromatography with atmospheric presure chemical ionization tandem mass spectrometry'
$ws.Range('G54').Value = 'vocab:1028'
$ws.Range('A55').Value = 'vocab:1031'
$ws.Range('B55').Value = 'gas chromatography coupled to electron capture negative ion mass spectrometry'
$ws.Range('G55').Value = 'vocab:1028'
$ws.Range('A56').Value = 'vocab:1032'
$ws.Range('B56').Value = 'gas chromatography with mass spectrometry with Negative chemical ionization'
$ws.Range('G56').Value = 'vocab:1028'
$ws.Range('A57').Value = 'vocab:1033'
$ws.Range('B57').Value = 'gas chromatography with high resolution mass spectrometry (HRMS)'
$ws.Range('G57').Value = 'vocab:1028'
$ws.Range('A58').Value = 'vocab:1034'
$ws.Range('B58').Value = 'gas chromatography with atmospheric presure chemical ionization high resolution mass spectrometry'
$ws.Range('G58').Value = 'vocab:1028'
$ws.Range('A59').Value = 'vocab:1035'
$ws.Range('B59').Value = 'gas chromatography with electron capture detector'
$ws.Range('G59').Value = 'vocab:1027'

Write-Output 'edit complete'
